$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55: fix politeness_score (column B) to be a proper number instead of text "3"
$ws.Range("B55").Value = 3

# Row 56: new annotation row from Ruilin
$ws.Range("A56").Value = "Ruilin"

# Column B on this row keeps the original (text) quirk, matching the source data
$ws.Range("B56").NumberFormat = "@"
$ws.Range("B56").Value = "3"
$ws.Range("B56").Style = "Normal"

$ws.Range("C56").Value = "无"
$ws.Range("D56").Value = "DIS"
$ws.Range("E56").Value = "RES"
$ws.Range("F56").Value = "e3eeb88f-0832-4aa9-a6cc-39ada0451b32"
$ws.Range("G56").Value = "BkJ3ibb0-_annotated.xlsx"
$ws.Range("H56").Value = "This paper shows that models trained on a synthetic dataset are vulnerable to small adversarial perturbations which lie on the data manifold."
